# "actualizacion estado del tp"
#
# - Row 8  (Iniciar Viaje Pendiente (Simulación)) / B8  : note -> "Probar al final"
# - Row 9  (Detener Viaje Pendiente (Simulación)) / B9  : note -> "Probar al final"
# - Row 12 (Cargar destinos...)                   / B12 : status green -> red
# - Row 13 (Excepciones)                          / B13 : status red -> green
# - Row 14 (Serializacion / XML)                  / B14 : status red -> green
# - Row 15 (Validaciones)                         / B15 : note "Faltan numeros negativos"
#                                                          -> "Faltan numeros negativos y patente"
# - Column B width widened (23.29 -> ~33.7 chars)
# - Active selection moved from B2 to B17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B so the longer note text fits.
$ws.Columns.Item(2).ColumnWidth = 32.8

# Update the note text on the two "pending" rows and the validations row.
$ws.Range("B8").Value = "Probar al final"
$ws.Range("B9").Value = "Probar al final"
$ws.Range("B15").Value = "Faltan numeros negativos y patente"

# Flip status colors: green = FF00B050, red = FFFF0000 (OLE BGR encoding).
$green = 0 + 176 * 256 + 80 * 65536
$red   = 255 + 0 * 256 + 0 * 65536

$ws.Range("B12").Interior.Color = $red
$ws.Range("B13").Interior.Color = $green
$ws.Range("B14").Interior.Color = $green

# Move the active selection to B17 (below the last data row).
[void]$ws.Range("B17").Select()
